$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new daily rows (2025-10-23, serial 45953) for the two stations,
# matching the formatting of the last existing row (45) by copying its
# formats down rather than re-deriving number formats (keeps styles.xml
# untouched / reuses the existing style indices).
$xlPasteFormats = -4122

$ws.Range("A45:F45").Copy()
$ws.Range("A46:F46").PasteSpecial($xlPasteFormats)

$ws.Range("A45:F45").Copy()
$ws.Range("A47:F47").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

$ws.Cells.Item(46, 1).Value = 45953
$ws.Cells.Item(46, 2).Value = "四方坪站"
$ws.Cells.Item(46, 3).Value = 10022.34
$ws.Cells.Item(46, 4).Value = 8006.63
$ws.Cells.Item(46, 5).Value = 3555.21
$ws.Cells.Item(46, 6).Value = 421

$ws.Cells.Item(47, 1).Value = 45953
$ws.Cells.Item(47, 2).Value = "高岭站"
$ws.Cells.Item(47, 3).Value = 5362.37
$ws.Cells.Item(47, 4).Value = 4575.08
$ws.Cells.Item(47, 5).Value = 1389.28
$ws.Cells.Item(47, 6).Value = 214

# Update the view: scroll so the new rows are visible and select M45 - this
# matches where the workbook was left selected/scrolled to after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M45").Select()
